$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated strikeout counts ("K" column, column G) regenerated from save_data source
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")
$kValues = @{
    2 = 1
    3 = 4
    4 = 2
    5 = 1
    6 = 1
    7 = 3
    8 = 3
    9 = 2
    10 = 0
    11 = 2
    12 = 0
    13 = 2
    14 = 0
    15 = 3
    16 = 2
    17 = 1
    18 = 1
    19 = 4
    20 = 4
    21 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 3
    32 = 2
    33 = 0
    34 = 0
    35 = 1
    36 = 2
    37 = 3
    38 = 1
    39 = 0
    40 = 2
    41 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 5
    46 = 1
    47 = 4
    48 = 2
    50 = 3
    51 = 0
    52 = 5
    53 = 5
    54 = 1
    55 = 4
    56 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

